# Word COM-interop script applying the tracked changes to
# "Intro Text Data Analytics presentation.docx"
#
# Summary of the edit:
#  - Append a reviewer comment (in blue, accent1-ish color) to the paragraph
#    that begins "We will analyse the impact of stringency ..." as two new
#    runs (Word often splits inserted text into multiple runs).
#  - Word re-drops a "_GoBack" bookmark at the last edited location when the
#    file is saved, so add that too (empty bookmark at the insertion point).

$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Locate the target paragraph ("We will analyse the impact ...").
# ------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.StartsWith("We will analyse the impact of stringency")) {
        $target = $para
        break
    }
}

# Position right after the existing text, but before the paragraph mark.
$insertPos = $target.Range.End - 1
$insertRange = $d.Range($insertPos, $insertPos)

# ------------------------------------------------------------------
# 2. Insert the first run of the reviewer comment.
# ------------------------------------------------------------------
$part1 = "(Der Satz ist mir nicht ganz klar: würd hier v"
$insertRange.InsertAfter($part1)
$run1Start = $insertPos
$run1End = $run1Start + $part1.Length
$run1Range = $d.Range($run1Start, $run1End)
$run1Range.Font.Name = "Calibri"
$run1Range.Font.Size = 11
$run1Range.Font.Color = 12874308

# ------------------------------------------------------------------
# 3. Insert the second run of the reviewer comment (continuation).
# ------------------------------------------------------------------
$part2 = "orschlagen, eher zu sagen, impact of stringency on death per million and happiness – dann könnten wir zum einen sagen, ob das harte Regime überhaupt etwas gebracht hat bzgl Anzahl Tote und ob es einen Effekt gab auf die Happiness und ggf sagen, ob es sich «rentiert» hat eine hohe Stringenfcy zu haben.)"
$run2Start = $run1End
$run2InsertPoint = $d.Range($run2Start, $run2Start)
$run2InsertPoint.InsertAfter($part2)
$run2End = $run2Start + $part2.Length
$run2Range = $d.Range($run2Start, $run2End)
$run2Range.Font.Name = "Calibri"
$run2Range.Font.Size = 11
$run2Range.Font.Color = 12874308

# ------------------------------------------------------------------
# 4. Drop the "_GoBack" bookmark at the end of the inserted text
#    (Word's automatic "last edit" bookmark).
# ------------------------------------------------------------------
$goBackRange = $d.Range($run2End, $run2End)
$d.Bookmarks.Add("_GoBack", $goBackRange)

Write-Host "Edit complete"
